# TC11_Canine_Filter_SamplePatho-Undefined.xlsx - "startup" sheet edit
#
# CasesTab / SamplesTab / FilesTab rows each carry a "query" (col B) and a
# "StatQuery" (col C). This updates:
#   * StatQuery (C2/C3/C4) for all three rows -> new Programs/Studies/Cases/
#     Samples/Case Files/Study Files count query (replaces the old
#     number_of_files/number_of_sample/number_of_cases/number_of_study one).
#   * CasesTab query (B2) -> adds a trailing `Cohort` column to the RETURN.
#   * FilesTab query (B4) -> drops the trailing `Study Code` column from the
#     RETURN (SamplesTab query in B3 is untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`r`nOPTIONAL MATCH (samp:sample)-->(c)`r`nOPTIONAL MATCH (diag:diagnosis)-->(c)`r`nOPTIONAL MATCH (f:file)-[*]->(c)`r`nOPTIONAL MATCH (sf:file)-->(s)`r`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`r`n WHERE samp.specific_sample_pathology IN [`"Undefined`"]  `r`nRETURN  `r`n    count(distinct p) AS Programs,`r`n    count(distinct s) AS Studies,`r`n    count(distinct c) AS Cases,`r`n    count(distinct samp) AS Samples,`r`n    count(distinct f) AS ``Case Files``,`r`n    count(distinct sf) AS ``Study Files``"
$ws.Cells.Item(2, 3).Value2 = $newStatQuery
$ws.Cells.Item(3, 3).Value2 = $newStatQuery
$ws.Cells.Item(4, 3).Value2 = $newStatQuery

$ws.Cells.Item(2, 2).Value2 = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`r`nMATCH (c)<--(diag:diagnosis)`r`nMATCH (samp:sample)-->(c) `r`n WHERE samp.specific_sample_pathology IN [`"Undefined`"]  `r`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`r`n  WITH DISTINCT c, s, demo, diag, co`r`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`r`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`r`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`r`n        coalesce(demo.breed, '') AS Breed ,`r`n        coalesce(diag.disease_term, '') AS Diagnosis ,`r`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`r`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`r`n        coalesce(demo.sex, '') AS Sex ,`r`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`r`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`r`n        coalesce(diag.best_response, '') AS ``Response to Treatment``,`r`ncoalesce(co.cohort_description, '') AS ``Cohort``"

$ws.Cells.Item(4, 2).Value2 = "MATCH (f:file)-->(parent)`r`nWITH DISTINCT f, parent`r`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`r`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`r`n MATCH (samp:sample)-->(c) `r`n WHERE samp.specific_sample_pathology IN [`"Undefined`"]  `r`nWITH DISTINCT f, parent, c, demo, diag, s`r`nRETURN coalesce(f.file_name, '') AS ``File Name``, `r`n        coalesce(labels(parent)[0], '') AS ``Association``,`r`n        coalesce(f.file_description, '') AS ``Description``,`r`n        coalesce(f.file_format, '') AS ``Format``,`r`n        coalesce(f.file_size, '') AS ``Size``,`r`n        coalesce(c.case_id, '') AS ``Case ID``, `r`n        coalesce(diag.disease_term,'') AS Diagnosis "

# Row heights grow/shrink to fit the new wrapped query text.
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# The author's save left the active cell on C2.
$ws.Range("C2").Select()
